$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 701 (shifts existing rows 701:770 down to 702:771)
$ws.Rows(701).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A701").Value = 6
$ws.Range("B701").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C701").Value = "Metropolitana"
$ws.Range("D701").Value = 45194
$ws.Range("E701").Value = 13
$ws.Range("F701").Value = 100112039
$ws.Range("G701").Value = "Ciboulette"
$ws.Range("H701").Value = "Sin especificar"
$ws.Range("I701").Value = "Primera"
$ws.Range("J701").Value = 530
$ws.Range("K701").Value = 900
$ws.Range("L701").Value = 1000
$ws.Range("M701").Value = 957
$ws.Range("N701").Value = "`$/docena de atados"
$ws.Range("O701").Value = "Región Metropolitana"
$ws.Range("P701").Value = 319
$ws.Range("Q701").Value = 3
$ws.Range("R701").Value = "Hortaliza"
